$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.757380843162537
$ws.Range("B1").Value = 2.58054518699646
$ws.Range("C1").Value = 2.827648401260376
$ws.Range("D1").Value = 3.394381523132324
$ws.Range("E1").Value = 1.206219911575317
